$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column M ("Дата отправки результатов" / ${experiment.sentDate}) - shifts former
# column N left into M, updates dimension/merged cells/used range automatically.
$ws.Columns("M").Delete()

# Update the jx comment annotations that reference the old last column "N" (now "M",
# since the sheet only goes up to column M after the delete). These are free-text
# comment bodies, so Excel does not auto-update them on column delete.
$ws.Range("A1").Comment.Text("Roman93:`njx:area(lastCell=`"M8`")")
$ws.Range("A3").Comment.Text("Roman93:`njx:if(condition=`"not empty(report.searchQuery)`", lastCell=`"M3`", areas=[`"A3:B3`"])")
$ws.Range("A4").Comment.Text("Roman93:`njx:each(items=`"report.filters`" var=`"filter`" lastCell=`"M4`")")
$ws.Range("A7").Comment.Text("Roman93:`njx:each(items=`"report.items`" var=`"experiment`" lastCell=`"M7`")")

# Move the active selection to A7 (matches the saved view state in the edited workbook).
$ws.Range("A7").Select()
